$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1068.421
$ws.Range("I4").Value = 581.2
$ws.Range("K4").Value = 581.2
$ws.Range("M4").Value = -467.2
$ws.Range("H8").Value = 55.57143
$ws.Range("I8").Value = 53
$ws.Range("J8").Value = 59
$ws.Range("K8").Value = 159
$ws.Range("L8").Value = 177
$ws.Range("M8").Value = -20
$ws.Range("N8").Value = -455
$ws.Range("H15").Value = 710.23254
$ws.Range("I15").Value = 710.23254
$ws.Range("K15").Value = 2130.69762
$ws.Range("M15").Value = -1961.69762
$ws.Range("H17").Value = 1637
$ws.Range("I17").Value = 1380
$ws.Range("J17").Value = 1650.5264
$ws.Range("K17").Value = 4140
$ws.Range("L17").Value = 4951.5792
$ws.Range("M17").Value = -3972
$ws.Range("N17").Value = -5287.5792
$ws.Range("H19").Value = 1214.1666
$ws.Range("J19").Value = 1369.7
$ws.Range("L19").Value = 1369.7
$ws.Range("N19").Value = -1719.7
$ws.Range("H51").Value = 7668.08
$ws.Range("I51").Value = 19100
$ws.Range("K51").Value = 19100
$ws.Range("M51").Value = -18616
$ws.Range("H58").Value = 4500.75
$ws.Range("I58").Value = 859
$ws.Range("J58").Value = 7333.222
$ws.Range("K58").Value = 2577
$ws.Range("L58").Value = 21999.666
$ws.Range("M58").Value = -2427
$ws.Range("N58").Value = -22299.666
$ws.Range("H62").Value = 5439.8096
$ws.Range("I62").Value = 5202.8667
$ws.Range("K62").Value = 5202.8667
$ws.Range("M62").Value = -4578.8667
$ws.Range("H64").Value = 6270.643
$ws.Range("I64").Value = 6245.6
$ws.Range("J64").Value = 6284.5557
$ws.Range("K64").Value = 6245.6
$ws.Range("L64").Value = 6284.5557
$ws.Range("M64").Value = -5997.6
$ws.Range("N64").Value = -6780.5557
$ws.Range("H65").Value = 5439.8096
$ws.Range("I65").Value = 5202.8667
$ws.Range("K65").Value = 26014.3335
$ws.Range("M65").Value = -22894.3335
$ws.Range("H67").Value = 6270.643
$ws.Range("I67").Value = 6245.6
$ws.Range("J67").Value = 6284.5557
$ws.Range("K67").Value = 6245.6
$ws.Range("L67").Value = 6284.5557
$ws.Range("M67").Value = -5387.6
$ws.Range("N67").Value = -8000.5557
$ws.Range("H74").Value = 3210.3333
$ws.Range("I74").Value = 3210.3333
$ws.Range("K74").Value = 3210.3333
$ws.Range("M74").Value = -2274.3333
$ws.Range("H77").Value = 3210.3333
$ws.Range("I77").Value = 3210.3333
$ws.Range("K77").Value = 16051.6665
$ws.Range("M77").Value = -11371.6665
$ws.Range("H86").Value = 12400
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 12400
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 12400
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -14646
$ws.Range("H89").Value = 12400
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 12400
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 62000
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -73232
$ws.Range("H92").Value = 3173.9092
$ws.Range("I92").Value = 1231
$ws.Range("K92").Value = 1231
$ws.Range("M92").Value = 17
$ws.Range("H94").Value = 3314.2222
$ws.Range("I94").Value = 3314.2222
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3314.2222
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2863.2222
$ws.Range("N94").Value = ""
$ws.Range("H106").Value = 5404.4165
$ws.Range("I106").Value = 5368.4546
$ws.Range("K106").Value = 5368.4546
$ws.Range("M106").Value = -4737.4546
$ws.Range("H107").Value = 989.41174
$ws.Range("I107").Value = 381.16666
$ws.Range("K107").Value = 381.16666
$ws.Range("M107").Value = 1538.83334
$ws.Range("H132").Value = 1872.341
$ws.Range("I132").Value = 1826.4147
$ws.Range("K132").Value = 5479.2441
$ws.Range("M132").Value = -2949.2441
$ws.Range("H137").Value = 2509.913
$ws.Range("I137").Value = 1498.4286
$ws.Range("K137").Value = 4495.2858
$ws.Range("M137").Value = -1945.2858
$ws.Range("H138").Value = 3117.8113
$ws.Range("I138").Value = 3043.625
$ws.Range("J138").Value = 3149.8918
$ws.Range("K138").Value = 9130.875
$ws.Range("L138").Value = 9449.6754
$ws.Range("M138").Value = -3990.875
$ws.Range("N138").Value = -19729.6754
$ws.Range("H141").Value = 4545.9443
$ws.Range("I141").Value = 4545.9443
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 13637.8329
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -8457.832900000001
$ws.Range("N141").Value = ""

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 657.2857
$ws.Range("I5").Value = 600.1667
$ws.Range("K5").Value = 600.1667
$ws.Range("M5").Value = -488.1667
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""
$ws.Range("H32").Value = 6610.068
$ws.Range("I32").Value = 6005.268
$ws.Range("K32").Value = 6005.268
$ws.Range("M32").Value = -5718.268
$ws.Range("H35").Value = 125012560
$ws.Range("I35").Value = 11737.667
$ws.Range("K35").Value = 11737.667
$ws.Range("M35").Value = -11331.667
$ws.Range("H37").Value = 500041820
$ws.Range("J37").Value = 500041820
$ws.Range("L37").Value = 500041820
$ws.Range("N37").Value = -500042366
$ws.Range("H45").Value = 2692.6667
$ws.Range("I45").Value = 1126.3334
$ws.Range("J45").Value = 4259
$ws.Range("K45").Value = 1126.3334
$ws.Range("L45").Value = 4259
$ws.Range("M45").Value = -749.3334
$ws.Range("N45").Value = -5013
$ws.Range("H74").Value = 2501.0715
$ws.Range("I74").Value = 1501.9
$ws.Range("J74").Value = 4999
$ws.Range("K74").Value = 1501.9
$ws.Range("L74").Value = 4999
$ws.Range("M74").Value = -627.9000000000001
$ws.Range("N74").Value = -6747
$ws.Range("H77").Value = 2501.0715
$ws.Range("I77").Value = 1501.9
$ws.Range("J77").Value = 4999
$ws.Range("K77").Value = 7509.5
$ws.Range("L77").Value = 24995
$ws.Range("M77").Value = -3141.5
$ws.Range("N77").Value = -33731
$ws.Range("H88").Value = 2194.6316
$ws.Range("I88").Value = 1442.2222
$ws.Range("J88").Value = 2871.8
$ws.Range("K88").Value = 1442.2222
$ws.Range("L88").Value = 2871.8
$ws.Range("M88").Value = -1036.2222
$ws.Range("N88").Value = -3683.8
$ws.Range("H91").Value = 2194.6316
$ws.Range("I91").Value = 1442.2222
$ws.Range("J91").Value = 2871.8
$ws.Range("K91").Value = 1442.2222
$ws.Range("L91").Value = 2871.8
$ws.Range("M91").Value = -38.22219999999993
$ws.Range("N91").Value = -5679.8
$ws.Range("H102").Value = 4108
$ws.Range("I102").Value = 2225.5833
$ws.Range("J102").Value = 9755.25
$ws.Range("K102").Value = 2225.5833
$ws.Range("L102").Value = 9755.25
$ws.Range("M102").Value = -603.5832999999998
$ws.Range("N102").Value = -12999.25
$ws.Range("H122").Value = 5399.5557
$ws.Range("I122").Value = 5637.125
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 16911.375
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -14461.375
$ws.Range("N122").Value = -15397
$ws.Range("H132").Value = 2274165
$ws.Range("I132").Value = 1443.8889
$ws.Range("J132").Value = 12501410
$ws.Range("K132").Value = 4331.6667
$ws.Range("L132").Value = 37504230
$ws.Range("M132").Value = -1801.6667
$ws.Range("N132").Value = -37509290

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 657.2857
$ws.Range("I4").Value = 600.1667
$ws.Range("K4").Value = 600.1667
$ws.Range("M4").Value = -485.1667
$ws.Range("H22").Value = 1383.4286
$ws.Range("I22").Value = 1180.8334
$ws.Range("K22").Value = 1180.8334
$ws.Range("M22").Value = -1007.8334
$ws.Range("H86").Value = 2195.4583
$ws.Range("I86").Value = 1565.0555
$ws.Range("J86").Value = 4086.6667
$ws.Range("K86").Value = 1565.0555
$ws.Range("L86").Value = 4086.6667
$ws.Range("M86").Value = -442.0554999999999
$ws.Range("N86").Value = -6332.6667
$ws.Range("H89").Value = 2195.4583
$ws.Range("I89").Value = 1565.0555
$ws.Range("J89").Value = 4086.6667
$ws.Range("K89").Value = 7825.2775
$ws.Range("L89").Value = 20433.3335
$ws.Range("M89").Value = -2209.2775
$ws.Range("N89").Value = -31665.3335
$ws.Range("H99").Value = 2489.3076
$ws.Range("I99").Value = 1823
$ws.Range("K99").Value = 1823
$ws.Range("M99").Value = -325
$ws.Range("H105").Value = 498512.2
$ws.Range("I105").Value = 678869.8
$ws.Range("K105").Value = 678869.8
$ws.Range("M105").Value = -677122.8
$ws.Range("H107").Value = 4660.6
$ws.Range("I107").Value = 5386.75
$ws.Range("K107").Value = 5386.75
$ws.Range("M107").Value = -3466.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3481.2727
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224
$ws.Range("H7").Value = 42.52381
$ws.Range("I7").Value = 50.6
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 50.6
$ws.Range("L7").Value = 40
$ws.Range("M7").Value = 62.4
$ws.Range("N7").Value = -266
$ws.Range("H31").Value = 32260946
$ws.Range("I31").Value = 40002330
$ws.Range("K31").Value = 40002330
$ws.Range("M31").Value = -40002035
$ws.Range("H34").Value = 32260946
$ws.Range("I34").Value = 40002330
$ws.Range("K34").Value = 40002330
$ws.Range("M34").Value = -40002128
$ws.Range("H99").Value = 18043.941
$ws.Range("J99").Value = 27118.875
$ws.Range("L99").Value = 27118.875
$ws.Range("N99").Value = -30114.875
$ws.Range("H122").Value = 2050.0833
$ws.Range("I122").Value = 2219.2222
$ws.Range("K122").Value = 6657.6666
$ws.Range("M122").Value = -4207.6666
$ws.Range("H126").Value = 18043.941
$ws.Range("J126").Value = 27118.875
$ws.Range("L126").Value = 81356.625
$ws.Range("N126").Value = -86296.625
$ws.Range("H134").Value = 1471.4375
$ws.Range("I134").Value = 1213.3214
$ws.Range("K134").Value = 3639.9642
$ws.Range("M134").Value = -1104.9642

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 22166
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("H38").Value = 600.3333
$ws.Range("I38").Value = 375
$ws.Range("J38").Value = 1051
$ws.Range("K38").Value = 1125
$ws.Range("L38").Value = 3153
$ws.Range("M38").Value = -778
$ws.Range("N38").Value = -3847
$ws.Range("H39").Value = 15465.8
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 6000
$ws.Range("M39").Value = -5706
$ws.Range("H40").Value = 186.42857
$ws.Range("I40").Value = 68
$ws.Range("J40").Value = 233.8
$ws.Range("K40").Value = 272
$ws.Range("L40").Value = 935.2
$ws.Range("M40").Value = -203
$ws.Range("N40").Value = -1073.2
$ws.Range("H42").Value = 166673580
$ws.Range("J42").Value = 9778
$ws.Range("L42").Value = 29334
$ws.Range("N42").Value = -30402
$ws.Range("H44").Value = 18107
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 18107
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 54321
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = -55117
$ws.Range("H57").Value = 33333
$ws.Range("J57").Value = 33333
$ws.Range("L57").Value = 99999
$ws.Range("N57").Value = -101117
$ws.Range("H69").Value = 13729.286
$ws.Range("I69").Value = 3160
$ws.Range("J69").Value = 21656.25
$ws.Range("K69").Value = 9480
$ws.Range("L69").Value = 64968.75
$ws.Range("M69").Value = -8669
$ws.Range("N69").Value = -66590.75
$ws.Range("H72").Value = 13729.286
$ws.Range("I72").Value = 3160
$ws.Range("J72").Value = 21656.25
$ws.Range("K72").Value = 28440
$ws.Range("L72").Value = 194906.25
$ws.Range("M72").Value = -24384
$ws.Range("N72").Value = -203018.25
$ws.Range("H136").Value = 6101.222
$ws.Range("I136").Value = 763.3333
$ws.Range("J136").Value = 16777
$ws.Range("K136").Value = 2289.9999
$ws.Range("L136").Value = 50331
$ws.Range("M136").Value = 2810.0001
$ws.Range("N136").Value = -60531

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 33331.668
$ws.Range("I5").Value = 49995
$ws.Range("K5").Value = 49995
$ws.Range("M5").Value = -49883
$ws.Range("H70").Value = 11247.58
$ws.Range("I70").Value = 10570.786
$ws.Range("J70").Value = 11510.777
$ws.Range("K70").Value = 10570.786
$ws.Range("L70").Value = 11510.777
$ws.Range("M70").Value = -10300.786
$ws.Range("N70").Value = -12050.777
$ws.Range("H73").Value = 11247.58
$ws.Range("I73").Value = 10570.786
$ws.Range("J73").Value = 11510.777
$ws.Range("K73").Value = 10570.786
$ws.Range("L73").Value = 11510.777
$ws.Range("M73").Value = -9634.786
$ws.Range("N73").Value = -13382.777

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14041.7
$ws.Range("I22").Value = 55750.5
$ws.Range("J22").Value = 3614.5
$ws.Range("K22").Value = 55750.5
$ws.Range("L22").Value = 3614.5
$ws.Range("M22").Value = -55455.5
$ws.Range("N22").Value = -4204.5
$ws.Range("H27").Value = 14041.7
$ws.Range("I27").Value = 55750.5
$ws.Range("J27").Value = 3614.5
$ws.Range("K27").Value = 55750.5
$ws.Range("L27").Value = 3614.5
$ws.Range("M27").Value = -55643.5
$ws.Range("N27").Value = -3828.5
$ws.Range("H40").Value = 10999
$ws.Range("I40").Value = 10999
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 10999
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -10863
$ws.Range("N40").Value = ""
$ws.Range("H46").Value = 1043.091
$ws.Range("J46").Value = 1997.25
$ws.Range("L46").Value = 1997.25
$ws.Range("N46").Value = -2373.25
$ws.Range("H55").Value = 934.1667
$ws.Range("I55").Value = 408.81818
$ws.Range("J55").Value = 1378.6923
$ws.Range("K55").Value = 408.81818
$ws.Range("L55").Value = 1378.6923
$ws.Range("M55").Value = -235.81818
$ws.Range("N55").Value = -1724.6923
$ws.Range("H68").Value = 2607748
$ws.Range("I68").Value = 5209683.5
$ws.Range("K68").Value = 5209683.5
$ws.Range("M68").Value = -5208934.5
$ws.Range("H71").Value = 2607748
$ws.Range("I71").Value = 5209683.5
$ws.Range("K71").Value = 26048417.5
$ws.Range("M71").Value = -26044673.5
$ws.Range("H93").Value = 2927449.5
$ws.Range("I93").Value = 3114.111
$ws.Range("J93").Value = 5559351
$ws.Range("K93").Value = 3114.111
$ws.Range("L93").Value = 5559351
$ws.Range("M93").Value = -1866.111
$ws.Range("N93").Value = -5561847
$ws.Range("H100").Value = 27810144
$ws.Range("I100").Value = 2963
$ws.Range("K100").Value = 2963
$ws.Range("M100").Value = -2422
$ws.Range("H122").Value = 4696.4346
$ws.Range("I122").Value = 3487.6487
$ws.Range("J122").Value = 9665.888999999999
$ws.Range("K122").Value = 10462.9461
$ws.Range("L122").Value = 28997.667
$ws.Range("M122").Value = -8012.946100000001
$ws.Range("N122").Value = -33897.667
$ws.Range("H136").Value = 4833.273
$ws.Range("I136").Value = 1896
$ws.Range("K136").Value = 5688
$ws.Range("M136").Value = -3138

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 33199.6
$ws.Range("J2").Value = 16000
$ws.Range("L2").Value = 16000
$ws.Range("N2").Value = -16224
$ws.Range("H96").Value = 12251.571
$ws.Range("I96").Value = 11152.4
$ws.Range("J96").Value = 14999.5
$ws.Range("K96").Value = 11152.4
$ws.Range("L96").Value = 14999.5
$ws.Range("M96").Value = -9779.4
$ws.Range("N96").Value = -17745.5
$ws.Range("H107").Value = 9300.333000000001
$ws.Range("I107").Value = 5933.6665
$ws.Range("K107").Value = 17800.9995
$ws.Range("M107").Value = -15880.9995
$ws.Range("H122").Value = 2134.353
$ws.Range("I122").Value = 1949
$ws.Range("J122").Value = 2399.1428
$ws.Range("K122").Value = 5847
$ws.Range("L122").Value = 7197.428400000001
$ws.Range("M122").Value = -3397
$ws.Range("N122").Value = -12097.4284
$ws.Range("H132").Value = 246624.05
$ws.Range("I132").Value = 2723.0588
$ws.Range("J132").Value = 1431286
$ws.Range("K132").Value = 8169.176399999999
$ws.Range("L132").Value = 4293858
$ws.Range("M132").Value = -5639.176399999999
$ws.Range("N132").Value = -4298918

